$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.550.80"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.729.42"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4805"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06177"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").Value = "1.730.60"
$ws.Range("E10").Value = "  -0.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07180"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6081"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.528"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "26.558.09"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006962"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("E20").Value = "  -2.02%  "

$ws.Range("D21").Value = "1.952.52"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.513"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.786"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.227"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.775"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.403"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.967"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08003"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.685"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04496"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6297"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9065"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.050"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.399"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.14%  "

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.495"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3886"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.042"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.58%  "

$ws.Range("E46").Value = "  -2.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05383"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.825"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3409"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
